# Mark specific vocabulary rows as "Processed" in column C.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(11,12,13,14,15)
for ($r = 37; $r -le 1027; $r += 10) {
    $rows += $r
}

foreach ($r in $rows) {
    $ws.Cells.Item($r, 3).Value = "Processed"
}
